$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế ngày LONG XUYÊN")

# Column D holds Notion's "last_edited_time" timestamp. Three distinct
# timestamp blocks get bumped forward (~11 minutes later):
#   rows  2-15 : 2024-08-03T03:17 -> 2024-08-03T03:28
#   rows 16-53 : 2024-08-03T03:18 -> 2024-08-03T03:29
#   rows 54-94 : 2024-08-03T03:19 -> 2024-08-03T03:30
$ws.Range("D2:D15").Value = "2024-08-03T03:28:00.000Z"
$ws.Range("D16:D53").Value = "2024-08-03T03:29:00.000Z"
$ws.Range("D54:D94").Value = "2024-08-03T03:30:00.000Z"
